$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2025-09-30 Tuesday" "2025-10-01 Wednesday"

Replace-Text "474×4=" "829×3="
Replace-Text "235×3=" "270×9="
Replace-Text "589×7=" "700×7="
Replace-Text "943×5=" "196×4="
Replace-Text "135×9=" "264×2="

Replace-Text "154×3=" "993×3="
Replace-Text "479×4=" "130×6="
Replace-Text "910×6=" "326×8="
Replace-Text "876×9=" "251×7="
Replace-Text "628×5=" "273×4="

Replace-Text "485×6=" "900×7="
Replace-Text "893×4=" "284×7="
Replace-Text "299×6=" "677×4="
Replace-Text "536×6=" "748×3="
Replace-Text "250×4=" "850×7="

Replace-Text "527×3=" "424×2="
Replace-Text "340×3=" "782×5="
Replace-Text "499×5=" "155×3="
Replace-Text "178×8=" "120×7="
Replace-Text "300×2=" "396×3="

Replace-Text "157×2=" "911×7="
Replace-Text "392×5=" "140×2="
Replace-Text "971×3=" "861×5="
Replace-Text "155×7=" "268×8="
Replace-Text "336×8=" "325×6="
